# The post that used to live at row 233 ("「自分の胸に聞きなさい。それが全て知っています」")
# was removed from the blog, so delete its entire row. Excel shifts every
# subsequent row up by one and the used range shrinks from A1:C312 to A1:C311.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(233).Delete()
